$d = $word.ActiveDocument

# 1. Change the meeting time suffix from " AM" to " PM".
#    Replace " AM" -> " P" first, then insert a trailing "M" run so the
#    text ends up split across two runs (" P" and "M"), matching the
#    target OOXML.
$d.Content.Find.Execute(" AM", $true, $false, $false, $false, $false, $true, 1, $false, " P", 2)

$r = $d.Content
$r.Find.Execute(" P", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("M")

# 2. Re-apply (no-op) replacements on the Goals bullet text so that Word
#    collapses the previously-split runs that share identical formatting
#    back into single runs.
$d.Content.Find.Execute("Understand and discuss the purpose of all sensors to be installed and its application to the challenge requirements.", $true, $false, $false, $false, $false, $true, 1, $false, "Understand and discuss the purpose of all sensors to be installed and its application to the challenge requirements.", 2)

$d.Content.Find.Execute("Design and draw build 1.0.1 with sensors.", $true, $false, $false, $false, $false, $true, 1, $false, "Design and draw build 1.0.1 with sensors.", 2)

$d.Content.Find.Execute(" Install sensors onto EV3_Vers_1.0.", $true, $false, $false, $false, $false, $true, 1, $false, " Install sensors onto EV3_Vers_1.0.", 2)

# 3. Remove the four trailing empty paragraphs at the end of the document
#    (after the bookmarked paragraph, right before the sectPr).
$count = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item($count - 3)
$lastPara = $d.Paragraphs.Item($count)
$start = $startPara.Range.Start
$end = $lastPara.Range.End
$delRange = $d.Range($start, $end)
$delRange.Delete()
